$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the scope/value cells on row 2
$ws.Range("B2").Value = "all"
$ws.Range("D2").Value = "bob"

# Update the selected range to A2:F2 with active cell A2
$ws.Range("A2:F2").Select()
